$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 43
$ws.Range("H43").Value = 5529.773
$ws.Range("J43").Value = 5697.3125
$ws.Range("L43").Value = 5697.3125
$ws.Range("N43").Value = -5835.3125

# Row 106
$ws.Range("H106").Value = 3484.45
$ws.Range("I106").Value = 3484.45
$ws.Range("K106").Value = 3484.45
$ws.Range("M106").Value = -2853.45

# Row 135
$ws.Range("H135").Value = 1349.8235
$ws.Range("I135").Value = 1428.0834
$ws.Range("J135").Value = 1162
$ws.Range("K135").Value = 12852.7506
$ws.Range("L135").Value = 10458
$ws.Range("M135").Value = -10317.7506
$ws.Range("N135").Value = -15528

# Row 137
$ws.Range("H137").Value = 1515.1936
$ws.Range("I137").Value = 1527.7142
$ws.Range("J137").Value = 1398.3334
$ws.Range("K137").Value = 4583.142599999999
$ws.Range("L137").Value = 4195.0002
$ws.Range("M137").Value = -2033.142599999999
$ws.Range("N137").Value = -9295.0002

# Row 141
$ws.Range("H141").Value = 4951.3076
$ws.Range("I141").Value = 4986.9
$ws.Range("J141").Value = 4832.6665
$ws.Range("K141").Value = 14960.7
$ws.Range("L141").Value = 14497.9995
$ws.Range("M141").Value = -9780.699999999999
$ws.Range("N141").Value = -24857.9995

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 3
$ws.Range("H3").Value = 1792.25
$ws.Range("I3").Value = 1792.25
$ws.Range("K3").Value = 1792.25
$ws.Range("M3").Value = -1677.25

# Row 32
$ws.Range("H32").Value = 3431.7542
$ws.Range("I32").Value = 2872.2034
$ws.Range("J32").Value = 19938.5
$ws.Range("K32").Value = 2872.2034
$ws.Range("L32").Value = 19938.5
$ws.Range("M32").Value = -2585.2034
$ws.Range("N32").Value = -20512.5

# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = ""

# Row 61
$ws.Range("H61").Value = 4260.6665
$ws.Range("I61").Value = 4123.75
$ws.Range("K61").Value = 4123.75
$ws.Range("M61").Value = -3911.75

# Row 74
$ws.Range("H74").Value = 2049.3684
$ws.Range("I74").Value = 2041
$ws.Range("J74").Value = 2200
$ws.Range("K74").Value = 2041
$ws.Range("L74").Value = 2200
$ws.Range("M74").Value = -1167
$ws.Range("N74").Value = -3948

# Row 77
$ws.Range("H77").Value = 2049.3684
$ws.Range("I77").Value = 2041
$ws.Range("J77").Value = 2200
$ws.Range("K77").Value = 10205
$ws.Range("L77").Value = 11000
$ws.Range("M77").Value = -5837
$ws.Range("N77").Value = -19736

# Row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = ""

# Row 132
$ws.Range("H132").Value = 2417.375
$ws.Range("I132").Value = 2417.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7252.125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4722.125
$ws.Range("N132").Value = ""

# Row 136
$ws.Range("H136").Value = 4260.6665
$ws.Range("I136").Value = 4123.75
$ws.Range("K136").Value = 12371.25
$ws.Range("M136").Value = -9821.25

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
# Row 5
$ws.Range("H5").Value = 1040.8
$ws.Range("J5").Value = 1300
$ws.Range("L5").Value = 1300
$ws.Range("N5").Value = -1526

# Row 7
$ws.Range("H7").Value = 2655.1667
$ws.Range("I7").Value = 1990
$ws.Range("J7").Value = 2788.2
$ws.Range("K7").Value = 1990
$ws.Range("L7").Value = 2788.2
$ws.Range("M7").Value = -1877
$ws.Range("N7").Value = -3014.2

# Row 8
$ws.Range("H8").Value = 983.3333
$ws.Range("I8").Value = 1225
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 1225
$ws.Range("L8").Value = 500
$ws.Range("M8").Value = -1085
$ws.Range("N8").Value = -780

# Row 17
$ws.Range("H17").Value = 5669.6665
$ws.Range("J17").Value = 5669.6665
$ws.Range("L17").Value = 5669.6665
$ws.Range("N17").Value = -6013.6665

# Row 94
$ws.Range("H94").Value = 7000
$ws.Range("I94").Value = 6000
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 6000
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -5549
$ws.Range("N94").Value = -10902

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 2
$ws.Range("H2").Value = 600
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -426

# Row 11
$ws.Range("H11").Value = 3283.3333
$ws.Range("J11").Value = 3283.3333
$ws.Range("L11").Value = 3283.3333
$ws.Range("N11").Value = -3563.3333

# Row 15
$ws.Range("H15").Value = 3625
$ws.Range("I15").Value = 2500
$ws.Range("J15").Value = 4000
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 4000
$ws.Range("M15").Value = -2330
$ws.Range("N15").Value = -4340

# Row 99
$ws.Range("H99").Value = 2821.7144
$ws.Range("J99").Value = 2000
$ws.Range("L99").Value = 2000
$ws.Range("N99").Value = -4996

# Row 105
$ws.Range("H105").Value = 3368.111
$ws.Range("I105").Value = 2664
$ws.Range("J105").Value = 5832.5
$ws.Range("K105").Value = 2664
$ws.Range("L105").Value = 5832.5
$ws.Range("M105").Value = -917
$ws.Range("N105").Value = -9326.5

# Row 126
$ws.Range("H126").Value = 2821.7144
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940

# Row 134
$ws.Range("H134").Value = 4172.0527
$ws.Range("I134").Value = 3184.7334
$ws.Range("J134").Value = 7874.5
$ws.Range("K134").Value = 9554.200199999999
$ws.Range("L134").Value = 23623.5
$ws.Range("M134").Value = -7019.200199999999
$ws.Range("N134").Value = -28693.5

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 18
$ws.Range("H18").Value = 5142.5713
$ws.Range("I18").Value = 1399
$ws.Range("K18").Value = 4197
$ws.Range("M18").Value = -4028

# Row 55
$ws.Range("H55").Value = 298130.94
$ws.Range("J55").Value = 4264.125
$ws.Range("L55").Value = 12792.375
$ws.Range("N55").Value = -13146.375

# Row 140
$ws.Range("H140").Value = 1581.1111
$ws.Range("I140").Value = 1581.1111
$ws.Range("K140").Value = 4743.3333
$ws.Range("M140").Value = 436.6666999999998

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 97
$ws.Range("H97").Value = 452.94116
$ws.Range("I97").Value = 427
$ws.Range("K97").Value = 427
$ws.Range("M97").Value = 69

# Row 126
$ws.Range("H126").Value = 3037.5
$ws.Range("I126").Value = 3037.5
$ws.Range("K126").Value = 9112.5
$ws.Range("M126").Value = -6642.5

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 87
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42246

# Row 90
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -131232

# Row 93
$ws.Range("H93").Value = 4580
$ws.Range("I93").Value = 2966.6667
$ws.Range("J93").Value = 7000
$ws.Range("K93").Value = 2966.6667
$ws.Range("L93").Value = 7000
$ws.Range("M93").Value = -1718.6667
$ws.Range("N93").Value = -9496

# Row 100
$ws.Range("H100").Value = 5336.0454
$ws.Range("I100").Value = 3221
$ws.Range("K100").Value = 3221
$ws.Range("M100").Value = -2680

# Row 136
$ws.Range("H136").Value = 34779.188
$ws.Range("I136").Value = 2196
$ws.Range("J136").Value = 49589.727
$ws.Range("K136").Value = 6588
$ws.Range("L136").Value = 148769.181
$ws.Range("M136").Value = -4038
$ws.Range("N136").Value = -153869.181

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 126
$ws.Range("H126").Value = 1838.7142
$ws.Range("I126").Value = 1329.2
$ws.Range("K126").Value = 3987.6
$ws.Range("M126").Value = -1517.6
